# "75% of Levels done"
# Edits the level-grid on Sheet1:
#  - clear two stray cells in the header row
#  - fix a couple of mis-typed "W" cells to blank space in row 2
#  - extend the "outline" marker column S down through rows 30-39
#  - fill in wall ("#") / door-key ("*") symbols for the newly-drawn
#    level rows 35-38
#  - add a brand new (blank-but-marked) row 39
#  - move the active-cell selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 1: remove the two stray blank-space cells -----------------------
$ws.Range("AY1").Clear()
$ws.Range("AZ1").Clear()

# --- Row 2: correct three cells (were "W", now blank space) --------------
$ws.Range("AX2").Value = " "
$ws.Range("AY2").Value = " "
$ws.Range("BA2").Value = " "

# --- Column S marker, rows 30-38 ------------------------------------------
$sCol = @(30, 31, 32, 33, 34, 35, 36, 37, 38)
foreach ($r in $sCol) {
    $ws.Range("S$r").Value = " "
}

# --- Row 35: wall/key symbols ---------------------------------------------
$ws.Range("F35").Value = "#"
$ws.Range("G35").Value = "#"
$ws.Range("K35").Value = "*"
$ws.Range("O35").Value = "#"
$ws.Range("P35").Value = "#"

# --- Row 36: wall symbols --------------------------------------------------
$ws.Range("F36").Value = "#"
$ws.Range("G36").Value = "#"
$ws.Range("H36").Value = "#"
$ws.Range("N36").Value = "#"
$ws.Range("O36").Value = "#"
$ws.Range("P36").Value = "#"

# --- Row 37: wall symbols (columns F-J, L-P) -------------------------------
$ws.Range("F37").Value = "#"
$ws.Range("G37").Value = "#"
$ws.Range("H37").Value = "#"
$ws.Range("I37").Value = "#"
$ws.Range("J37").Value = "#"
$ws.Range("L37").Value = "#"
$ws.Range("M37").Value = "#"
$ws.Range("N37").Value = "#"
$ws.Range("O37").Value = "#"
$ws.Range("P37").Value = "#"

# --- Row 38 (new row): wall symbols (columns F-P) --------------------------
$ws.Range("F38").Value = "#"
$ws.Range("G38").Value = "#"
$ws.Range("H38").Value = "#"
$ws.Range("I38").Value = "#"
$ws.Range("J38").Value = "#"
$ws.Range("K38").Value = "#"
$ws.Range("L38").Value = "#"
$ws.Range("M38").Value = "#"
$ws.Range("N38").Value = "#"
$ws.Range("O38").Value = "#"
$ws.Range("P38").Value = "#"

# --- Row 39 (brand new row): blank-space markers D-S ------------------------
$row39cols = @("D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S")
foreach ($c in $row39cols) {
    $ws.Range("$c`39").Value = " "
}

# --- Move the active selection --------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("BI18").Select() | Out-Null
